$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Coinranking price/volume columns are stored as plain text
    # (e.g. "56.884.04", "0.0982", "  +0.90%  "). Assigning .Value
    # directly lets Excel auto-detect a numeric type for anything
    # that parses as a number (rounding/reformatting it), so force
    # the cell to Text first, then restore the original "Normal"
    # style so no stray number-format style sticks to the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "56.884.04"
Set-TextValue "E2" "  +0.90%  "
Set-TextValue "D3" "2.396.73"
Set-TextValue "E3" "  +0.79%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "507.06"
Set-TextValue "E5" "  +1.62%  "
Set-TextValue "D6" "134.60"
Set-TextValue "E6" "  +4.33%  "
Set-TextValue "D7" "0.998"
Set-TextValue "D8" "0.554"
Set-TextValue "E8" "  +0.54%  "
Set-TextValue "D9" "2.404.34"
Set-TextValue "E9" "  -0.03%  "
Set-TextValue "D10" "0.0982"
Set-TextValue "E10" "  +2.64%  "
Set-TextValue "D11" "0.151"
Set-TextValue "E11" "  +0.37%  "
Set-TextValue "D12" "0.338"
Set-TextValue "E12" "  +6.46%  "
Set-TextValue "D13" "4.68"
Set-TextValue "E13" "  +0.31%  "
Set-TextValue "D14" "2.823.54"
Set-TextValue "E14" "  +0.60%  "
Set-TextValue "D15" "56.870.23"
Set-TextValue "E15" "  +1.10%  "
Set-TextValue "D16" "21.89"
Set-TextValue "E16" "  +1.83%  "
Set-TextValue "D17" "0.0000134"
Set-TextValue "E17" "  +2.19%  "
Set-TextValue "D18" "2.368.77"
Set-TextValue "E18" "  -1.54%  "
Set-TextValue "D19" "10.21"
Set-TextValue "E19" "  +1.14%  "
Set-TextValue "D20" "4.07"
Set-TextValue "E20" "  +1.18%  "
Set-TextValue "D21" "311.17"
Set-TextValue "E21" "  +0.56%  "
Set-TextValue "E22" "  -0.17%  "
Set-TextValue "E23" "  +0.18%  "
Set-TextValue "D24" "5.65"
Set-TextValue "E24" "  +1.69%  "
Set-TextValue "D25" "65.71"
Set-TextValue "E25" "  +1.19%  "
Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  +0.07%  "
Set-TextValue "D27" "0.375"
Set-TextValue "E27" "  -0.02%  "
Set-TextValue "D28" "0.152"
Set-TextValue "E28" "  +1.16%  "
Set-TextValue "D29" "7.37"
Set-TextValue "E29" "  +1.47%  "
Set-TextValue "D30" "173.37"
Set-TextValue "E30" "  +0.59%  "
Set-TextValue "D31" "0.0₃0733"
Set-TextValue "E31" "  +2.86%  "
Set-TextValue "D32" "1.66"
Set-TextValue "E32" "  +0.63%  "
Set-TextValue "E33" "  +2.70%  "
Set-TextValue "D34" "5.88"
Set-TextValue "E34" "  -3.96%  "
Set-TextValue "E35" "  +0.15%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  +0.22%  "
Set-TextValue "D37" "17.95"
Set-TextValue "E37" "  +0.63%  "
Set-TextValue "D38" "1.21"
Set-TextValue "E38" "  -0.66%  "
Set-TextValue "D39" "3.88"
Set-TextValue "E39" "  +2.47%  "
Set-TextValue "D40" "36.67"
Set-TextValue "E40" "  +2.24%  "
Set-TextValue "D41" "0.821"
Set-TextValue "E41" "  +4.23%  "
Set-TextValue "E42" "  +1.46%  "
Set-TextValue "D43" "133.04"
Set-TextValue "E43" "  +2.77%  "
Set-TextValue "D44" "3.41"
Set-TextValue "E44" "  +2.41%  "
Set-TextValue "D45" "5.04"
Set-TextValue "E45" "  +3.99%  "
Set-TextValue "D46" "0.569"
Set-TextValue "E46" "  +0.83%  "
Set-TextValue "D47" "0.0911"
Set-TextValue "E47" "  +1.46%  "
Set-TextValue "D48" "249.81"
Set-TextValue "E48" "  -0.84%  "
Set-TextValue "D49" "0.0487"
Set-TextValue "E49" "  +0.70%  "
Set-TextValue "E50" "  +2.05%  "
Set-TextValue "D51" "17.34"
Set-TextValue "E51" "  +7.50%  "
